$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text-format on the Price column cells that are being updated so that
# numeric-looking strings (e.g. "675.03", "35.90") are preserved verbatim
# as text instead of being auto-coerced to floating point numbers. The
# NumberFormat is cleared again immediately after the write so the cell
# keeps its original (default) style.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.658.11'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.57%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.809.21'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.23%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '675.03'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +8.19%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.11'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.03%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.807.44'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.24%  '

$ws.Range('E8').Value = '  -0.05%  '

$ws.Range('E9').Value = '  +0.92%  '

$ws.Range('E10').Value = '  +1.56%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.26'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +7.66%  '

$ws.Range('E12').Value = '  +0.41%  '

$ws.Range('E13').Value = '  -0.57%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.90'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.02%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.452.51'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.04%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.811.73'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.74%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.710.01'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.57%  '

$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.64'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.11%  '

$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.20'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.42%  '

$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.114'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.45%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.31'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +18.57%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '478.11'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.53%  '

$ws.Range('E23').Value = '  +1.17%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.41'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.73%  '

$ws.Range('E25').Value = '  -3.16%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.22'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.69%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.34'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.35%  '

$ws.Range('E28').Value = '  -1.88%  '

$ws.Range('E29').Value = '  +0.07%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.960.58'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.20%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.91'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +9.18%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.29'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.70%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.38'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.38%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.51'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.86%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.180'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.47%  '

$ws.Range('E36').Value = '  +2.08%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.18%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.765.57'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.33%  '

$ws.Range('E39').Value = '  +1.05%  '

$ws.Range('E40').Value = '  +1.00%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.93'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.07%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.967'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.05%  '

$ws.Range('E43').Value = '  -0.08%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.12'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +12.01%  '

$ws.Range('E45').Value = '  -0.01%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '45.96'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +6.76%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '159.42'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.90%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '48.17'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.26%  '

$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.43'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.52%  '

$ws.Range('B50').Value = 'TheGraph'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.300'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.69%  '

$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000289'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +6.24%  '
